$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Dashboard" worksheet between "Admin" and "PIM".
#    We copy "Admin" (rather than Worksheets.Add()) so the new sheet inherits
#    the same namespace / sheetFormatPr plumbing the other sheets already
#    carry, then we wipe its contents and rebuild them from scratch.
#
#    NOTE: inserting a sheet shifts tab indices, and this host's Worksheet
#    handles are index-bound -- any handle obtained *before* the insert can
#    silently resolve to the wrong tab afterwards. So we fetch handles again
#    immediately after the structural change and use only those from here on.
# ---------------------------------------------------------------------------
$adminSheet0 = $wb.Worksheets.Item("Admin")
$adminSheet0.Copy($null, $adminSheet0)
$dashSheet0 = $wb.Worksheets.Item("Admin (2)")
$dashSheet0.Name = "Dashboard"

$adminSheet = $wb.Worksheets.Item("Admin")
$dashSheet  = $wb.Worksheets.Item("Dashboard")
$pimSheet   = $wb.Worksheets.Item("PIM")

$dashSheet.Cells.Clear()

# ---------------------------------------------------------------------------
# 2. Populate the Dashboard sheet.
#    Column order of writes matters: it reproduces the shared-string table
#    order seen in the target workbook (the author filled most of row 2,
#    circled back to add a 12th tab column, then filled in two columns that
#    had been skipped on the first pass).
#
#    Values are always written *before* PasteSpecial formatting is applied --
#    setting .Value after a format-only paste can make this host recompute
#    (and drop) an outdated quotePrefix style.
# ---------------------------------------------------------------------------

# Row 1 header cells A1..L1 (style copied from Admin!A1 / Admin!B1)
$dashSheet.Range("A1").Value = "TC01_Dashboard_verifyDasboardMainMenus"
$adminSheet.Range("A1").Copy()
$dashSheet.Range("A1").PasteSpecial(-4122)

$dashSheet.Range("B1").Value = "TabName1"
$dashSheet.Range("C1").Value = "TabName2"
$dashSheet.Range("D1").Value = "TabName3"
$dashSheet.Range("E1").Value = "TabName4"
$dashSheet.Range("F1").Value = "TabName5"
$dashSheet.Range("G1").Value = "TabName6"
$dashSheet.Range("H1").Value = "TabName7"
$dashSheet.Range("I1").Value = "TabName8"
$dashSheet.Range("J1").Value = "TabName9"
$dashSheet.Range("K1").Value = "TabName10"
$dashSheet.Range("L1").Value = "TabName11"
$adminSheet.Range("B1").Copy()
$dashSheet.Range("B1:L1").PasteSpecial(-4122)

# Row 2 data cells
$dashSheet.Range("B2").Value = "Admin"
$adminSheet.Range("B2").Copy()
$dashSheet.Range("B2").PasteSpecial(-4122)

$dashSheet.Range("C2").Value = "PIM"
$dashSheet.Range("D2").Value = "Leave"
$dashSheet.Range("F2").Value = "Recruitment"
$dashSheet.Range("G2").Value = "My Info"
$dashSheet.Range("I2").Value = "Dashboard"
$dashSheet.Range("J2").Value = "Directory"
$dashSheet.Range("K2").Value = "Maintenance"
$dashSheet.Range("L2").Value = "Claim"
$dashSheet.Range("M2").Value = "Buzz"

# F2 ("Recruitment") carries the same style as Admin!E2 ("ErrorMessage").
$adminSheet.Range("E2").Copy()
$dashSheet.Range("F2").PasteSpecial(-4122)
$dashSheet.Range("F2").Value = "Recruitment"

# 12th tab column, added after the rest of row 2 was filled in.
$dashSheet.Range("M1").Value = "TabName12"
$adminSheet.Range("B1").Copy()
$dashSheet.Range("M1").PasteSpecial(-4122)

$dashSheet.Range("E2").Value = "Time"
$dashSheet.Range("H2").Value = "Performance"

$dashSheet.Range("N2").Value = "TC01_Dashboard_verifyDasboardMainMenus"
$adminSheet.Range("A1").Copy()
$dashSheet.Range("N2").PasteSpecial(-4122)

# Column A width + selection.
$dashSheet.Columns.Item(2).ColumnWidth = $dashSheet.StandardWidth
$dashSheet.Columns.Item(3).ColumnWidth = $dashSheet.StandardWidth
$dashSheet.Columns.Item(4).ColumnWidth = $dashSheet.StandardWidth
$dashSheet.Columns.Item(1).ColumnWidth = 29.45
$dashSheet.Range("K12").Select()

# ---------------------------------------------------------------------------
# 3. Add two new rows of test data to the PIM sheet.
# ---------------------------------------------------------------------------
$pimSheet.Range("A5").Value = "TC03_PIM_CreateEmployeeWithLoginDetails"
$pimSheet.Range("A3").Copy()
$pimSheet.Range("A5").PasteSpecial(-4122)

$pimSheet.Range("B5").Value = "EmployeeLastName"
$pimSheet.Range("C5").Value = "NewEmployeePassword"
$pimSheet.Range("B3").Copy()
$pimSheet.Range("B5:C5").PasteSpecial(-4122)

$pimSheet.Range("B6").Value = "Kapoor"
$pimSheet.Range("C6").Value = "R29vZHdpbGwxMjM0NQ=="

$pimSheet.Range("D6").Value = "TC03_PIM_CreateEmployeeWithLoginDetails"
$pimSheet.Range("C4").Copy()
$pimSheet.Range("D6").PasteSpecial(-4122)

$pimSheet.Range("D14").Select()

# ---------------------------------------------------------------------------
# 4. Selection / active-tab bookkeeping on the Admin sheet, then re-activate
#    PIM so it becomes the workbook's active tab (matches the target diff).
# ---------------------------------------------------------------------------
$adminSheet.Range("A1:G2").Select()

$pimSheet.Activate()
